$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.690.14"
$ws.Range("E2").Value = "  -1.89%  "

$ws.Range("D3").Value = "3.167.79"
$ws.Range("E3").Value = "  -4.87%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.57%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.15%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "3.164.98"
$ws.Range("E8").Value = "  -4.93%  "

$ws.Range("E9").Value = "  -0.67%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.140"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.84%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.23"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.21%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.452"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.81%  "

$ws.Range("E13").Value = "  -5.27%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.79"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.17%  "

$ws.Range("D15").Value = "3.694.36"
$ws.Range("E15").Value = "  -4.90%  "

$ws.Range("E16").Value = "  -1.17%  "

$ws.Range("D17").Value = "3.170.43"
$ws.Range("E17").Value = "  -5.00%  "

$ws.Range("D18").Value = "62.710.28"
$ws.Range("E18").Value = "  -2.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "459.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.84%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.693"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.47%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.62"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.68%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.92%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.46%  "

$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("E27").Value = "  +0.01%  "

$ws.Range("E28").Value = "  -4.40%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.27%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.43%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.02"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.55%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.34%  "

$ws.Range("E33").Value = "  -4.31%  "

$ws.Range("E34").Value = "  -7.14%  "

$ws.Range("E35").Value = "  -7.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.80"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.86%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.20"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.45%  "

$ws.Range("D38").Value = "0.0₃0701"
$ws.Range("E38").Value = "  -7.03%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0387"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.45%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "402.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.52%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.05"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.05%  "

$ws.Range("E42").Value = "  -5.97%  "

$ws.Range("E43").Value = "  -5.78%  "

$ws.Range("D44").Value = "2.788.18"
$ws.Range("E44").Value = "  -11.10%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.250"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.87%  "

$ws.Range("E47").Value = "  -7.45%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.21"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.44%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.15"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.21%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.36"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.15%  "

$ws.Range("E51").Value = "  -2.38%  "
